$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 34, shifting existing rows 34-44 down to 35-45
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new weekly price record
$ws.Range("A34").Value = 4
$ws.Range("B34").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C34").Value = "Los Lagos"
$ws.Range("D34").Value = 44627
$ws.Range("E34").Value = 10
$ws.Range("F34").Value = 100112043
$ws.Range("G34").Value = "Pepino dulce"
$ws.Range("H34").Value = "Cultivar IV Región"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 50
$ws.Range("K34").Value = 18000
$ws.Range("L34").Value = 18000
$ws.Range("M34").Value = 18000
$ws.Range("N34").Value = "$/bandeja 18 kilos"
$ws.Range("O34").Value = "Provincia de Limarí"
$ws.Range("P34").Value = 1000
$ws.Range("Q34").Value = 18
$ws.Range("R34").Value = "Hortaliza"
